# Weekly Forward Look stats update 24.10.25
#
# 1. Bump the "as at <date>" reference in the intro line from
#    17 October 2025 to 24 October 2025.
# 2. Remove the row for "Criminal justice statistics quarterly: June 2025"
#    (row 5) now that it has been published - deleting the entire row
#    shifts every subsequent row up by one, which also drops the final
#    (now-empty) week row at the bottom of the table.
# 3. The conditional-formatting rules covering the table are anchored to
#    explicit ranges (A5:F63 / A5:A63) that don't auto-shrink when a row
#    is deleted, so re-point them at the new, one-row-shorter extents
#    (A5:F62 / A5:A62).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "This list contains a week-by-week view of  MoJ Official and National Statistics that have been pre-announced on the gov.uk release calendar as at 24 October 2025"

$ws.Rows(5).Delete()

$fcs = $ws.Range("A5:F62").FormatConditions
for ($i = 1; $i -le $fcs.Count; $i++) {
    $fc = $fcs.Item($i)
    $addr = $fc.AppliesTo.Address()
    if ($addr -eq '$A$5:$F$63') {
        $fc.ModifyAppliesToRange($ws.Range("A5:F62"))
    } elseif ($addr -eq '$A$5:$A$63') {
        $fc.ModifyAppliesToRange($ws.Range("A5:A62"))
    }
}
